$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.338815331459045
$ws.Range("B1").Value = 1.398788213729858
$ws.Range("C1").Value = 3.883192777633667
$ws.Range("D1").Value = 3.321572780609131
$ws.Range("E1").Value = 1.065452098846436
